$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("K2").Value = 59.6
$ws.Range("N2").Value = 54.82400714602223

# Row 3
$ws.Range("D3").Value = 92070.92999999999
$ws.Range("E3").Value = 61.3
$ws.Range("F3").Value = 1.34
$ws.Range("K3").Value = 55.6
$ws.Range("N3").Value = 54.82400714602223

# Row 4
$ws.Range("K4").Value = 51.4
$ws.Range("N4").Value = 54.82400714602223

# Row 5
$ws.Range("K5").Value = 49.6
$ws.Range("N5").Value = 54.82400714602223

# Row 6
$ws.Range("K6").Value = 35.8
$ws.Range("N6").Value = 54.82400714602223
